$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '26.053.61'
Set-TextValue $ws.Range("E2") '  +0.01%  '
Set-TextValue $ws.Range("D3") '1.648.58'
Set-TextValue $ws.Range("E3") '  +0.41%  '
Set-TextValue $ws.Range("E4") '  -0.23%  '
Set-TextValue $ws.Range("D5") '218.10'
Set-TextValue $ws.Range("E5") '  +0.29%  '
Set-TextValue $ws.Range("D6") '0.5198'
Set-TextValue $ws.Range("E6") '  +0.36%  '
Set-TextValue $ws.Range("E7") '  -0.24%  '
Set-TextValue $ws.Range("E8") '  +0.69%  '
Set-TextValue $ws.Range("E9") '  +0.63%  '
Set-TextValue $ws.Range("D10") '20.34'
Set-TextValue $ws.Range("E10") '  -0.15%  '
Set-TextValue $ws.Range("D11") '0.07653'
Set-TextValue $ws.Range("E11") '  -1.38%  '
Set-TextValue $ws.Range("D12") '4.577'
Set-TextValue $ws.Range("E12") '  +2.50%  '
Set-TextValue $ws.Range("D13") '1.658.30'
Set-TextValue $ws.Range("E13") '  +1.05%  '
Set-TextValue $ws.Range("D14") '1.875.77'
Set-TextValue $ws.Range("E14") '  +0.42%  '
Set-TextValue $ws.Range("D15") '0.5577'
Set-TextValue $ws.Range("E15") '  +0.44%  '
Set-TextValue $ws.Range("D16") '0.0₅8103'
Set-TextValue $ws.Range("E16") '  +1.71%  '
Set-TextValue $ws.Range("D17") '65.12'
Set-TextValue $ws.Range("E17") '  +0.75%  '
Set-TextValue $ws.Range("D18") '26.045.83'
Set-TextValue $ws.Range("E18") '  +0.01%  '
Set-TextValue $ws.Range("E19") '  -0.25%  '
Set-TextValue $ws.Range("D20") '4.604'
Set-TextValue $ws.Range("E20") '  -0.13%  '
Set-TextValue $ws.Range("D21") '10.46'
Set-TextValue $ws.Range("E21") '  +3.98%  '
Set-TextValue $ws.Range("D22") '190.93'
Set-TextValue $ws.Range("E22") '  -0.83%  '
Set-TextValue $ws.Range("E23") '  -0.56%  '
Set-TextValue $ws.Range("E24") '  -0.25%  '
Set-TextValue $ws.Range("D25") '144.22'
Set-TextValue $ws.Range("E25") '  -1.47%  '
Set-TextValue $ws.Range("E26") '  -1.62%  '
Set-TextValue $ws.Range("D27") '7.174'
Set-TextValue $ws.Range("E27") '  +0.37%  '
Set-TextValue $ws.Range("D28") '15.84'
Set-TextValue $ws.Range("E28") '  -0.10%  '
Set-TextValue $ws.Range("D29") '1.505'
Set-TextValue $ws.Range("E29") '  +1.68%  '
Set-TextValue $ws.Range("D30") '0.05350'
Set-TextValue $ws.Range("D31") '1.268'
Set-TextValue $ws.Range("E31") '  +0.35%  '
Set-TextValue $ws.Range("D32") '3.457'
Set-TextValue $ws.Range("E32") '  -0.10%  '
Set-TextValue $ws.Range("E33") '  -0.52%  '
Set-TextValue $ws.Range("D34") '1.548'
Set-TextValue $ws.Range("E34") '  -2.57%  '
Set-TextValue $ws.Range("D35") '2.419'
Set-TextValue $ws.Range("E35") '  +0.41%  '
Set-TextValue $ws.Range("D36") '2.779'
Set-TextValue $ws.Range("E36") '  -0.23%  '
Set-TextValue $ws.Range("D37") '0.9423'
Set-TextValue $ws.Range("E37") '  +0.72%  '
Set-TextValue $ws.Range("D38") '0.5614'
Set-TextValue $ws.Range("E38") '  -0.41%  '
Set-TextValue $ws.Range("D39") '0.01572'
Set-TextValue $ws.Range("E39") '  +0.03%  '
Set-TextValue $ws.Range("D40") '5.865'
Set-TextValue $ws.Range("E40") '  -1.20%  '
Set-TextValue $ws.Range("E41") '  -0.18%  '
Set-TextValue $ws.Range("D42") '1.030.05'
Set-TextValue $ws.Range("E42") '  -2.19%  '
Set-TextValue $ws.Range("D43") '0.8241'
Set-TextValue $ws.Range("E43") '  -1.67%  '
Set-TextValue $ws.Range("E44") '  -1.79%  '
Set-TextValue $ws.Range("D45") '1.788.27'
Set-TextValue $ws.Range("E45") '  +0.51%  '
Set-TextValue $ws.Range("D46") '0.0₈112'
Set-TextValue $ws.Range("E46") '  +5.74%  '
Set-TextValue $ws.Range("D47") '57.12'
Set-TextValue $ws.Range("E47") '  +0.54%  '
Set-TextValue $ws.Range("E48") '  -0.53%  '
Set-TextValue $ws.Range("E49") '  -0.33%  '
Set-TextValue $ws.Range("D50") '7.908'
Set-TextValue $ws.Range("E50") '  -0.45%  '
Set-TextValue $ws.Range("D51") '0.05141'
Set-TextValue $ws.Range("E51") '  -3.75%  '
